$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R2").Value = "ftp.sra.ebi.ac.uk/vol1/run/ERR147/ERR14789831/PAW77640.tar.gz"
$ws.Range("R3").Value = "ftp.sra.ebi.ac.uk/vol1/run/ERR147/ERR14789831/PAW77640.tar.gz"
$ws.Range("R4").Value = "ftp.sra.ebi.ac.uk/vol1/run/ERR147/ERR14789831/PAW77640.tar.gz"
$ws.Range("R5").Value = "ftp.sra.ebi.ac.uk/vol1/run/ERR147/ERR14789832/PAY12289_barcode13.tar.gz"
$ws.Range("R6").Value = "ftp.sra.ebi.ac.uk/vol1/run/ERR147/ERR14789833/PAS01578.tar.gz"
$ws.Range("R7").Value = "ftp.sra.ebi.ac.uk/vol1/run/ERR147/ERR14789833/PAS01578.tar.gz"
$ws.Range("R8").Value = "ftp.sra.ebi.ac.uk/vol1/run/ERR147/ERR14789834/PAY12289_barcode12.tar.gz"
$ws.Range("R9").Value = "ftp.sra.ebi.ac.uk/vol1/run/ERR147/ERR14789835/PAW78174_barcode01.tar.gz"
$ws.Range("R10").Value = "ftp.sra.ebi.ac.uk/vol1/run/ERR147/ERR14789836/PAW78174_barcode11.tar.gz"
$ws.Range("R11").Value = "ftp.sra.ebi.ac.uk/vol1/run/ERR147/ERR14789837/PAW78174_barcode05.tar.gz"
$ws.Range("R12").Value = "ftp.sra.ebi.ac.uk/vol1/run/ERR147/ERR14789838/PAW78174_barcode07.tar.gz"
$ws.Range("R13").Value = "ftp.sra.ebi.ac.uk/vol1/run/ERR147/ERR14789839/PAW78174_barcode06.tar.gz"
$ws.Range("R14").Value = "ftp.sra.ebi.ac.uk/vol1/run/ERR147/ERR14789840/PAW78174_barcode10.tar.gz"
$ws.Range("R15").Value = "ftp.sra.ebi.ac.uk/vol1/run/ERR147/ERR14789841/PAW78174_barcode08.tar.gz"
$ws.Range("R17").Value = "ftp.sra.ebi.ac.uk/vol1/run/ERR147/ERR14789842/PAW78174_barcode04.tar.gz"
$ws.Range("R18").Value = "ftp.sra.ebi.ac.uk/vol1/run/ERR147/ERR14789843/PAW78174_barcode03.tar.gz"
$ws.Range("R19").Value = "ftp.sra.ebi.ac.uk/vol1/run/ERR147/ERR14789844/PAW78174_barcode02.tar.gz"
